$wb = $excel.ActiveWorkbook

# Sheet 1: "Trends Status"
$ws1 = $wb.Worksheets.Item("Trends Status")
$ws1.Range("B2").Value = 98
$ws1.Range("C2").Value = 64
$ws1.Range("B3").Value = 106
$ws1.Range("C3").Value = 78
$ws1.Range("B4").Value = 423
$ws1.Range("C4").Value = 303
$ws1.Range("B5").Value = 185
$ws1.Range("C5").Value = 284
$ws1.Range("B6").Value = 98
$ws1.Range("C6").Value = 189
$ws1.Range("C7").Value = 17
$ws1.Range("C8").Value = 11

# Sheet 2: "Range Status"
$ws2 = $wb.Worksheets.Item("Range Status")
$ws2.Range("B3").Value = 70
$ws2.Range("B4").Value = 235

# Sheet 3: "Priority Status"
$ws3 = $wb.Worksheets.Item("Priority Status")
$ws3.Range("B2").Value = 199
$ws3.Range("B3").Value = 343
$ws3.Range("B4").Value = 404

# Sheet 4: "Species qualification"
$ws4 = $wb.Worksheets.Item("Species qualification")
$ws4.Range("B3").Value = 530
$ws4.Range("B4").Value = 650

# Sheet 5: "SoIB-IUCN cross-tab"
$ws5 = $wb.Worksheets.Item("SoIB-IUCN cross-tab")
$ws5.Range("B5").Value = 23
$ws5.Range("C5").Value = 33
$ws5.Range("B6").Value = 101
$ws5.Range("C6").Value = 299
$ws5.Range("D6").Value = 387
$ws5.Range("B7").Value = 1
$ws5.Range("C7").Value = 4
$ws5.Range("B8").Value = 199
$ws5.Range("C8").Value = 404
$ws5.Range("D8").Value = 343
